$wb = $excel.ActiveWorkbook

$todo = $wb.Worksheets.Item("Todo ")
$defect = $wb.Worksheets.Item("Defect")

# --- "Defect" sheet: update row 2 (HP bar defect) with PIC, finish date, status ---
$defect.Range("E2").Value = "Fish"
$defect.Range("H2").Value = "Solved"

$defect.Range("F2").Copy()
$defect.Range("G2").PasteSpecial(-4122)
$defect.Range("G2").Value = 43103

$defect.Range("C24").Select() | Out-Null

# --- "Todo " sheet: add a new Todo item in row 13 ---
$todo.Range("A13").Value = "Optimize player jump"
$todo.Range("C13").Value = "NA"
$todo.Range("D13").Value = "Fish"

$todo.Range("F12").Copy()
$todo.Range("F13").PasteSpecial(-4122)
$todo.Range("F13").Value = 43103

$todo.Range("D30").Select() | Out-Null
